# Generate Report for handoff
# The ce536fa2-b3da-4f14-a74d-08faf91cb8e1 file is now ready for handoff
# (previously it was "Handed back: in sync with en-US") for both the
# zh-cn and de-de locales. Update the Status + Latest Handoff Datetime
# columns on each locale sheet, and the rolled-up status on the Overview
# sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the ce536fa2-...md entry (zh-cn status in column
# B, de-de status in column C).
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: row 3 is the ce536fa2-...md entry.
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-01-25 10:56:38"

# de-de sheet: row 3 is the ce536fa2-...md entry.
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-01-25 10:56:48"
